$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the Heading1 title.
#    Build it by copying the formatted run structure (empty run + bold run)
#    from the duplicate bold paragraph near the end of the doc, so the
#    leading empty <w:r/> and the bold run formatting match exactly, then
#    retarget its text and append the non-bold remainder.
# ---------------------------------------------------------------------------
$title = $d.Paragraphs.First
$title.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$paraCount = $d.Paragraphs.Count
$boldDupPara = $d.Paragraphs.Item($paraCount - 1)
$metaPara.Range.FormattedText = $boldDupPara.Range.FormattedText

$metaRange = $metaPara.Range
$metaRange.MoveEnd(1, -1)
$metaRange.Text = "Meta description"

$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range
$metaRange.MoveEnd(1, -1)
$metaRange.InsertAfter(": Read our review of 100 Fortunes, an expanding slot game with a unique Chinese aesthetic. Play for free and enjoy bonus features like free spins and expanding Wild symbols.")

# ---------------------------------------------------------------------------
# 2. Remove the now-duplicate bold "Play 100 Fortunes..." paragraph that used
#    to sit right before the closing italic paragraph.
# ---------------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$boldDupPara = $d.Paragraphs.Item($paraCount - 1)
$boldDupPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3. Replace the text of the closing italic paragraph with the new image
#    generation prompt, keeping the italic formatting. Use the Range.Text
#    setter (not Find/Replace, which smart-quotes straight quotes) and
#    shrink the range off the paragraph mark first to avoid the end-of-story
#    insert/duplicate quirk.
# ---------------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$italicPara = $d.Paragraphs.Item($paraCount)
$italicRange = $italicPara.Range
$italicRange.MoveEnd(1, -1)
$italicRange.Text = 'Prompt: Create a feature image for "100 Fortunes" in a cartoon style featuring a happy Maya warrior with glasses. The image should depict the warrior standing in front of the game grid with a big smile on his face, holding a golden yo-yo and surrounded by Chinese cultural elements such as red paper lanterns and Ming vases. The cartoon style should be colorful and eye-catching, with a clear focus on the Maya warrior''s joyful expression. The image should convey the excitement and fun of playing this slot game while also highlighting the unique blend of Chinese and Maya cultures.'
